$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell $ws "D2" '67.085.46'
Set-TextCell $ws "E2" '  -3.35%  '

# Row 3
Set-TextCell $ws "D3" '3.515.70'
Set-TextCell $ws "E3" '  -4.25%  '

# Row 4
Set-TextCell $ws "D4" '1.00'
Set-TextCell $ws "E4" '  +0.17%  '

# Row 5
Set-TextCell $ws "D5" '608.93'
Set-TextCell $ws "E5" '  -5.76%  '

# Row 6
Set-TextCell $ws "D6" '153.65'
Set-TextCell $ws "E6" '  -3.98%  '

# Row 7
Set-TextCell $ws "D7" '3.511.88'
Set-TextCell $ws "E7" '  -4.26%  '

# Row 8
Set-TextCell $ws "E8" '  +0.06%  '

# Row 9
Set-TextCell $ws "D9" '0.485'
Set-TextCell $ws "E9" '  -2.37%  '

# Row 10
Set-TextCell $ws "E10" '  -2.98%  '

# Row 11
Set-TextCell $ws "D11" '6.87'
Set-TextCell $ws "E11" '  -3.22%  '

# Row 12
Set-TextCell $ws "D12" '0.428'
Set-TextCell $ws "E12" '  -3.12%  '

# Row 13
Set-TextCell $ws "D13" '0.0000221'
Set-TextCell $ws "E13" '  -4.86%  '

# Row 14
Set-TextCell $ws "B14" 'Avalanche'
Set-TextCell $ws "C14" 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell $ws "D14" '32.02'
Set-TextCell $ws "E14" '  -2.16%  '

# Row 15
Set-TextCell $ws "B15" 'WrappedliquidstakedEther2.0'
Set-TextCell $ws "C15" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell $ws "D15" '4.110.37'
Set-TextCell $ws "E15" '  -4.20%  '

# Row 16
Set-TextCell $ws "D16" '3.528.41'
Set-TextCell $ws "E16" '  -3.74%  '

# Row 17
Set-TextCell $ws "D17" '67.231.93'
Set-TextCell $ws "E17" '  -3.11%  '

# Row 18
Set-TextCell $ws "E18" '  +0.86%  '

# Row 19
Set-TextCell $ws "D19" '6.36'
Set-TextCell $ws "E19" '  -1.92%  '

# Row 20
Set-TextCell $ws "D20" '15.45'
Set-TextCell $ws "E20" '  -3.66%  '

# Row 21
Set-TextCell $ws "D21" '452.51'
Set-TextCell $ws "E21" '  -2.96%  '

# Row 22
Set-TextCell $ws "D22" '9.37'
Set-TextCell $ws "E22" '  -5.12%  '

# Row 23
Set-TextCell $ws "E23" '  -0.70%  '

# Row 24
Set-TextCell $ws "D24" '78.86'
Set-TextCell $ws "E24" '  -0.78%  '

# Row 25
Set-TextCell $ws "B25" 'Dai'
Set-TextCell $ws "C25" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws "D25" '1.00'
Set-TextCell $ws "E25" '  -0.21%  '

# Row 26
Set-TextCell $ws "B26" 'WrappedeETH'
Set-TextCell $ws "C26" 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextCell $ws "D26" '3.662.76'
Set-TextCell $ws "E26" '  -4.05%  '

# Row 27
Set-TextCell $ws "D27" '0.0000122'
Set-TextCell $ws "E27" '  -3.20%  '

# Row 28
Set-TextCell $ws "D28" '10.42'
Set-TextCell $ws "E28" '  -4.23%  '

# Row 29
Set-TextCell $ws "D29" '8.28'
Set-TextCell $ws "E29" '  -8.79%  '

# Row 30
Set-TextCell $ws "E30" '  -3.08%  '

# Row 31
Set-TextCell $ws "D31" '1.66'
Set-TextCell $ws "E31" '  -3.13%  '

# Row 32
Set-TextCell $ws "E32" '  +0.10%  '

# Row 33
Set-TextCell $ws "D33" '25.90'
Set-TextCell $ws "E33" '  -3.03%  '

# Row 34
Set-TextCell $ws "E34" '  -5.68%  '

# Row 35
Set-TextCell $ws "D35" '6.18'
Set-TextCell $ws "E35" '  -4.34%  '

# Row 36
Set-TextCell $ws "D36" '0.157'
Set-TextCell $ws "E36" '  -3.81%  '

# Row 37
Set-TextCell $ws "D37" '3.514.47'
Set-TextCell $ws "E37" '  -4.01%  '

# Row 38
Set-TextCell $ws "D38" '7.98'
Set-TextCell $ws "E38" '  -5.51%  '

# Row 39
Set-TextCell $ws "E39" '  -0.04%  '

# Row 40
Set-TextCell $ws "D40" '1.00'
Set-TextCell $ws "E40" '  +0.21%  '

# Row 41
Set-TextCell $ws "D41" '177.70'
Set-TextCell $ws "E41" '  -0.14%  '

# Row 42
Set-TextCell $ws "D42" '5.59'
Set-TextCell $ws "E42" '  -5.50%  '

# Row 43
Set-TextCell $ws "D43" '0.0878'
Set-TextCell $ws "E43" '  -2.94%  '

# Row 44
Set-TextCell $ws "D44" '2.11'
Set-TextCell $ws "E44" '  -3.23%  '

# Row 45
Set-TextCell $ws "D45" '0.890'
Set-TextCell $ws "E45" '  -3.79%  '

# Row 46
Set-TextCell $ws "D46" '29.24'
Set-TextCell $ws "E46" '  +7.71%  '

# Row 47
Set-TextCell $ws "D47" '45.64'
Set-TextCell $ws "E47" '  -2.11%  '

# Row 48
Set-TextCell $ws "E48" '  -3.52%  '

# Row 49
Set-TextCell $ws "B49" 'Cosmos'
Set-TextCell $ws "C49" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws "D49" '7.64'
Set-TextCell $ws "E49" '  -2.58%  '

# Row 50
Set-TextCell $ws "B50" 'ONDO'
Set-TextCell $ws "C50" 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextCell $ws "D50" '1.21'
Set-TextCell $ws "E50" '  -4.07%  '

# Row 51
Set-TextCell $ws "E51" '  -4.47%  '
